# Auto-generated scheduled-runner market data refresh for Anima_Profits sheets
$wb = $excel.ActiveWorkbook

$ALC = $wb.Worksheets.Item("ALC")
$ARM = $wb.Worksheets.Item("ARM")
$BSM = $wb.Worksheets.Item("BSM")
$CRP = $wb.Worksheets.Item("CRP")
$CUL = $wb.Worksheets.Item("CUL")
$GSM = $wb.Worksheets.Item("GSM")
$LTW = $wb.Worksheets.Item("LTW")

$ALC.Range("H43").Value = 1300.3478
$ALC.Range("I43").Value = 550.1111
$ALC.Range("J43").Value = 1782.6428
$ALC.Range("K43").Value = 550.1111
$ALC.Range("L43").Value = 1782.6428
$ALC.Range("M43").Value = -481.1111
$ALC.Range("N43").Value = -1920.6428
$ALC.Range("H137").Value = 1540.7727
$ALC.Range("I137").Value = 1308.1538
$ALC.Range("J137").Value = 1876.7778
$ALC.Range("K137").Value = 3924.4614
$ALC.Range("L137").Value = 5630.3334
$ALC.Range("M137").Value = -1374.4614
$ALC.Range("N137").Value = -10730.3334
$ALC.Range("H138").Value = 3033.889
$ALC.Range("I138").Value = 16598.5
$ALC.Range("J138").Value = 2402.9768
$ALC.Range("K138").Value = 49795.5
$ALC.Range("L138").Value = 7208.930399999999
$ALC.Range("M138").Value = -44655.5
$ALC.Range("N138").Value = -17488.9304
$ARM.Range("H43").Value = 9900
$ARM.Range("J43").Value = 9900
$ARM.Range("L43").Value = 9900
$ARM.Range("N43").Value = -10526
$ARM.Range("H123").Value = 33286
$ARM.Range("J123").Value = 33286
$ARM.Range("L123").Value = 33286
$ARM.Range("N123").Value = -43086
$ARM.Range("H134").Value = 78619.336
$ARM.Range("J134").Value = 78619.336
$ARM.Range("L134").Value = 78619.336
$ARM.Range("N134").Value = -88759.336
$ARM.Range("H140").Value = 50429
$ARM.Range("J140").Value = 50429
$ARM.Range("L140").Value = 50429
$ARM.Range("N140").Value = -60789
$ARM.Range("H141").Value = 100429
$ARM.Range("J141").Value = 100429
$ARM.Range("L141").Value = 100429
$ARM.Range("N141").Value = -110789
$BSM.Range("H134").Value = 3196.484
$BSM.Range("I134").Value = 3146.524
$BSM.Range("J134").Value = 3301.4
$BSM.Range("K134").Value = 9439.572
$BSM.Range("L134").Value = 9904.200000000001
$BSM.Range("M134").Value = -6904.572
$BSM.Range("N134").Value = -14974.2
$BSM.Range("H140").Value = 99280
$BSM.Range("J140").Value = 99280
$BSM.Range("L140").Value = 99280
$BSM.Range("N140").Value = -109640
$CRP.Range("H31").Value = 3041.7727
$CRP.Range("I31").Value = 1062.8125
$CRP.Range("J31").Value = 4904.3237
$CRP.Range("K31").Value = 1062.8125
$CRP.Range("L31").Value = 4904.3237
$CRP.Range("M31").Value = -767.8125
$CRP.Range("N31").Value = -5494.3237
$CRP.Range("H34").Value = 3041.7727
$CRP.Range("I34").Value = 1062.8125
$CRP.Range("J34").Value = 4904.3237
$CRP.Range("K34").Value = 1062.8125
$CRP.Range("L34").Value = 4904.3237
$CRP.Range("M34").Value = -860.8125
$CRP.Range("N34").Value = -5308.3237
$CRP.Range("H140").Value = 59600
$CRP.Range("J140").Value = 59600
$CRP.Range("L140").Value = 59600
$CRP.Range("N140").Value = -69960
$CUL.Range("H2").Value = 139.16667
$CUL.Range("I2").Value = 22.5
$CUL.Range("J2").Value = 197.5
$CUL.Range("K2").Value = 135
$CUL.Range("L2").Value = 1185
$CUL.Range("M2").Value = -22
$CUL.Range("N2").Value = -1411
$CUL.Range("H7").Value = 484.85715
$CUL.Range("I7").Value = 192.25
$CUL.Range("J7").Value = 875
$CUL.Range("K7").Value = 576.75
$CUL.Range("L7").Value = 2625
$CUL.Range("M7").Value = -464.75
$CUL.Range("N7").Value = -2849
$CUL.Range("H33").Value = 62612.438
$CUL.Range("I33").Value = 12624.875
$CUL.Range("K33").Value = 75749.25
$CUL.Range("M33").Value = -75466.25
$CUL.Range("H35").Value = 4798.25
$CUL.Range("J35").Value = 5440.857
$CUL.Range("L35").Value = 16322.571
$CUL.Range("N35").Value = -16898.571
$CUL.Range("H38").Value = 1133.25
$CUL.Range("I38").Value = 62.25
$CUL.Range("J38").Value = 1668.75
$CUL.Range("K38").Value = 186.75
$CUL.Range("L38").Value = 5006.25
$CUL.Range("M38").Value = 160.25
$CUL.Range("N38").Value = -5700.25
$CUL.Range("H40").Value = 236.45454
$CUL.Range("I40").Value = 250.1
$CUL.Range("J40").Value = 100
$CUL.Range("K40").Value = 1000.4
$CUL.Range("L40").Value = 400
$CUL.Range("M40").Value = -931.4
$CUL.Range("N40").Value = -538
$CUL.Range("H107").Value = 1786.1333
$CUL.Range("J107").Value = 2315.8728
$CUL.Range("L107").Value = 6947.6184
$CUL.Range("N107").Value = -10787.6184
$CUL.Range("H122").Value = 6351
$CUL.Range("I122").Value = 374.53845
$CUL.Range("J122").Value = 25774.5
$CUL.Range("K122").Value = 3370.84605
$CUL.Range("L122").Value = 231970.5
$CUL.Range("M122").Value = -920.8460500000001
$CUL.Range("N122").Value = -236870.5
$GSM.Range("H104").Value = 33000
$GSM.Range("J104").Value = 33000
$GSM.Range("L104").Value = 33000
$GSM.Range("N104").Value = -39988
$LTW.Range("H127").Value = 30000
$LTW.Range("J127").Value = 30000
$LTW.Range("L127").Value = 30000
$LTW.Range("N127").Value = -39920
$LTW.Range("H139").Value = 9680000
$LTW.Range("J139").Value = 60000
$LTW.Range("L139").Value = 60000
$LTW.Range("N139").Value = -70280
$LTW.Range("H141").Value = 80315
$LTW.Range("J141").Value = 80315
$LTW.Range("L141").Value = 80315
$LTW.Range("N141").Value = -90675
$ARM.Range("H39").Value = 15000
$ARM.Range("I39").Value = 0
$ARM.Range("J39").Value = 15000
$ARM.Range("K39").Value = 0
$ARM.Range("L39").Value = 15000
$ARM.Range("M39").ClearContents()
$ARM.Range("N39").Value = -16040
$ARM.Range("H109").Value = 55500
$ARM.Range("J109").Value = 55500
$ARM.Range("L109").Value = 55500
$ARM.Range("N109").Value = -58274
$ARM.Range("H125").Value = 100715
$ARM.Range("J125").Value = 100715
$ARM.Range("L125").Value = 100715
$ARM.Range("N125").Value = -110555
$BSM.Range("H59").Value = 0
$BSM.Range("J59").Value = 0
$BSM.Range("L59").Value = 0
$BSM.Range("N59").ClearContents()
$CRP.Range("H137").Value = 60226.668
$CRP.Range("I137").Value = 0
$CRP.Range("J137").Value = 60226.668
$CRP.Range("K137").Value = 0
$CRP.Range("L137").Value = 60226.668
$CRP.Range("M137").ClearContents()
$CRP.Range("N137").Value = -70426.66800000001

Write-Host "Applied Anima_Profits market data refresh"